$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.269.03"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "2.098.10"
$ws.Range("E3").Value = "  +4.40%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'251.05"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'54.47"
$ws.Range("E8").Value = "  +20.88%  "
$ws.Range("D9").Value = "'61.78"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").Value = "'0.375"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("E11").Value = "  +4.14%  "
$ws.Range("E12").Value = "  +7.37%  "
$ws.Range("D13").Value = "'15.35"
$ws.Range("E13").Value = "  +5.42%  "
$ws.Range("E14").Value = "  +4.36%  "
$ws.Range("D15").Value = "'0.841"
$ws.Range("E15").Value = "  +3.50%  "
$ws.Range("D16").Value = "2.101.10"
$ws.Range("E16").Value = "  +4.60%  "
$ws.Range("D17").Value = "'5.15"
$ws.Range("E17").Value = "  +5.08%  "
$ws.Range("D18").Value = "37.293.04"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("D19").Value = "'72.59"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").Value = "'14.55"
$ws.Range("E20").Value = "  +13.25%  "
$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  +2.85%  "
$ws.Range("D22").Value = "'241.69"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").Value = "'5.23"
$ws.Range("E23").Value = "  +7.05%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("D26").Value = "'171.37"
$ws.Range("E26").Value = "  +4.25%  "
$ws.Range("D27").Value = "'9.31"
$ws.Range("E27").Value = "  +8.38%  "
$ws.Range("D28").Value = "'20.68"
$ws.Range("E28").Value = "  +5.21%  "
$ws.Range("E29").Value = "  +4.12%  "
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("E31").Value = "  +26.71%  "
$ws.Range("D32").Value = "'22.36"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("D34").Value = "'0.0615"
$ws.Range("E34").Value = "  +4.91%  "
$ws.Range("D35").Value = "'0.0904"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.85"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.29"
$ws.Range("E38").Value = "  +5.29%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'4.11"
$ws.Range("E39").Value = "  +3.17%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "'18.53"
$ws.Range("E41").Value = "  +15.90%  "
$ws.Range("E42").Value = "  +3.79%  "
$ws.Range("D43").Value = "'1.16"
$ws.Range("E43").Value = "  +4.92%  "
$ws.Range("D44").Value = "'98.90"
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D45").Value = "'0.0920"
$ws.Range("E45").Value = "  +12.45%  "
$ws.Range("D46").Value = "'2.80"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'4.17"
$ws.Range("E47").Value = "  +107.26%  "
$ws.Range("D48").Value = "1.320.29"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").Value = "'2.97"
$ws.Range("E49").Value = "  +7.65%  "
$ws.Range("E50").Value = "  +14.85%  "
$ws.Range("D51").Value = "2.292.79"
$ws.Range("E51").Value = "  +4.50%  "
